$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Move to location (...)" text (strip trailing newline, change coordinates)
# and adjust the corresponding Robot (column B) boolean flags per the diff.

$ws.Range("A2").Value = "Move to location (9, 5) and remove the toolkit."

$ws.Range("A3").Value = "Move to location (3, 8) and remove the liquid spill."

$ws.Range("A4").Value = "Move to location (1, 4) and remove the large debris."
$ws.Range("B4").Value = $true

$ws.Range("A5").Value = "Move to location (6, 5) and remove the dust."
$ws.Range("B5").Value = $false

$ws.Range("A6").Value = "Move to location (9, 5) and remove the grass."

$ws.Range("A7").Value = "Move to location (5, 12) and remove the small debris."
$ws.Range("B7").Value = $false

$ws.Range("A8").Value = "Move to location (11, 12) and remove the vehicle."
$ws.Range("B8").Value = $true

$ws.Range("A9").Value = "Move to location (12, 1) and remove the construction materials."

$ws.Range("A10").Value = "Move to location (8, 12) and remove the tree branches."

$ws.Range("A11").Value = "Move to location (3, 5) and remove the screws."

$wb.Save()
